$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Led greenfield project building end-to-end scalable distributed machine learning platform using PyTorch, Ray, and Kubernetes for cancer detection from deep sequencing (methylated DNA) and protein data, enabling training of much larger models leveraging data distributed parallel (DDP) processing.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Led greenfield project building end-to-end scalable distributed machine learning platform using PyTorch, Ray, and Kubernetes for cancer detection from deep sequencing (methylated DNA) and protein data, enabling training of much larger models leveraging data distributed parallel (DDP) processing speeding up model training by >10x.",
    2
)

$d.Content.Find.Execute(
    "Deployed and managed an organization-wide MLFlow based model tracking system using Terraform, Pulumi, and Google Cloud enabling live monitoring of deep learning model training progress, instantaneous sharing results, and completely automated and reproducible report generation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Deployed and managed an organization-wide MLFlow based model tracking system using Terraform, Pulumi, and Google Cloud enabling live-monitoring of deep learning model training progress, instantaneous results sharing, and completely automated and reproducible report generation - reducing researcher manual effort by at least 5x.",
    2
)

$d.Content.Find.Execute(
    "Built scalable multitask learning, elastic net, and neural network based models in PyTorch with improved performance for classifying Colorectal Cancer risk from cell-free DNA data.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Built scalable multitask learning, elastic net, and neural network based models in PyTorch with improved performance for classifying Colorectal Cancer risk from cell-free DNA data for a clinical trial cohort of >27,000 individuals.",
    2
)

$d.Content.Find.Execute(
    "Automated performance metric report generation for all polygenic risk score classifiers leveraging MLFlow artifact storage and headless jupyter execution, reducing researcher time spent on analysis from days to minutes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Automated performance metric report generation for all polygenic risk score classifiers leveraging MLFlow artifact storage and headless Jupyter execution, reducing researcher time spent on analysis from days to minutes.",
    2
)
